$wb = $excel.ActiveWorkbook

# --- Sheet "Weekly Sales": append row 14 (new week: 2025-01-05) ---
$ws1 = $wb.Worksheets.Item("Weekly Sales")
$dateFmt1 = $ws1.Cells.Item(13, 1).NumberFormat

$ws1.Cells.Item(14, 1).Value = 45662.99999999999
$ws1.Cells.Item(14, 1).NumberFormat = $dateFmt1
$ws1.Cells.Item(14, 2).Value = 5

# --- Sheet "Merged (Optional)": append row 22 (new week: 2025-01-05) ---
$ws3 = $wb.Worksheets.Item("Merged (Optional)")
$dateFmt3 = $ws3.Cells.Item(21, 1).NumberFormat

$ws3.Cells.Item(22, 1).Value = 45662.99999999999
$ws3.Cells.Item(22, 1).NumberFormat = $dateFmt3
$ws3.Cells.Item(22, 2).Value = 5
$ws3.Cells.Item(22, 3).Value = 0
